$wb = $excel.ActiveWorkbook

# Rename sheets
$wb.Worksheets.Item("Letters_LTE_3_5oz").Name = "USPS_Letter_Rates"
$wb.Worksheets.Item("Flats_GT_QuarterInch").Name = "USPS_Flat_Rates"

# Activate USPS_Letter_Rates (was Letters_LTE_3_5oz) -> becomes the selected tab
$ws = $wb.Worksheets.Item("USPS_Letter_Rates")
$ws.Activate()

# Update the selection on that sheet to E6
$ws.Range("E6").Select()
